$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I0 / IF headers with the same style as the other header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: I0 / IF values for rows 2-78
$data = @{
    2 = @(5, 6)
    3 = @(9, 9)
    4 = @(9, 9)
    5 = @(8, 8)
    6 = @(6, 6)
    7 = @(7, 7)
    8 = @(6, 6)
    9 = @(5, 5)
    10 = @(6, 6)
    11 = @(5, 5)
    12 = @(8, 9)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(9, 9)
    16 = @(6, 7)
    17 = @(7, 7)
    18 = @(7, 7)
    19 = @(6, 7)
    20 = @(9, 10)
    21 = @(7, 7)
    22 = @(7, 7)
    23 = @(5, 6)
    24 = @(6, 7)
    25 = @(8, 8)
    26 = @(6, 6)
    27 = @(8, 8)
    28 = @(6, 6)
    29 = @(7, 7)
    30 = @(9, 9)
    31 = @(6, 6)
    32 = @(6, 7)
    33 = @(7, 7)
    34 = @(7, 7)
    35 = @(6, 6)
    36 = @(7, 7)
    37 = @(5, 6)
    38 = @(7, 8)
    39 = @(8, 8)
    40 = @(7, 8)
    41 = @(5, 5)
    42 = @(7, 7)
    43 = @(7, 7)
    44 = @(7, 7)
    45 = @(6, 6)
    46 = @(7, 7)
    47 = @(6, 6)
    48 = @(4, 5)
    49 = @(6, 7)
    50 = @(5, 6)
    51 = @(5, 6)
    52 = @(8, 8)
    53 = @(7, 7)
    54 = @(8, 8)
    55 = @(8, 8)
    56 = @(8, 8)
    57 = @(8, 8)
    58 = @(8, 8)
    59 = @(9, 9)
    60 = @(9, 9)
    61 = @(9, 9)
    62 = @(9, 9)
    63 = @(8, 8)
    64 = @(9, 9)
    65 = @(10, 10)
    66 = @(9, 9)
    67 = @(9, 9)
    68 = @(9, 9)
    69 = @(9, 9)
    70 = @(8, 9)
    71 = @(10, 10)
    72 = @(9, 9)
    73 = @(9, 9)
    74 = @(6, 6)
    75 = @(8, 8)
    76 = @(9, 9)
    77 = @(9, 9)
    78 = @(4, 4)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
